# Re-generate the blue gradient in the "Colour Code" column (B) of the
# "Date Colours" table so that it forms one smooth, continuous scale from
# row 2 through row 48. Previously rows 4-26 were out of sequence (an
# earlier edit had scrambled them) and rows 45-48 were missing their
# colour codes entirely - this restores/extends the gradient end to end.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colourCodes = @{
    2  = "#fcfcff"
    3  = "#f8faff"
    4  = "#f5f7ff"
    5  = "#f2f5ff"
    6  = "#eef2ff"
    7  = "#ebf0ff"
    8  = "#e7edff"
    9  = "#e4eaff"
    10 = "#e1e8ff"
    11 = "#dde5ff"
    12 = "#dae3ff"
    13 = "#d6e0ff"
    14 = "#d2deff"
    15 = "#cfdbff"
    16 = "#cbd9ff"
    17 = "#c8d6ff"
    18 = "#c4d4ff"
    19 = "#c0d1ff"
    20 = "#bdcfff"
    21 = "#b9ccff"
    22 = "#b5caff"
    23 = "#b1c7ff"
    24 = "#adc5ff"
    25 = "#a9c3ff"
    26 = "#a5c0ff"
    27 = "#a1beff"
    28 = "#9dbbff"
    29 = "#99b9ff"
    30 = "#95b6ff"
    31 = "#90b4ff"
    32 = "#8cb2ff"
    33 = "#87afff"
    34 = "#82adff"
    35 = "#7eabff"
    36 = "#79a8ff"
    37 = "#74a6ff"
    38 = "#6ea4ff"
    39 = "#69a1ff"
    40 = "#639fff"
    41 = "#5c9dff"
    42 = "#569aff"
    43 = "#4f98ff"
    44 = "#4796ff"
    45 = "#3e94ff"
    46 = "#3391ff"
    47 = "#268fff"
    48 = "#118dff"
}

foreach ($row in $colourCodes.Keys) {
    $ws.Cells.Item($row, 2).Value = $colourCodes[$row]
}
